$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 167, shifting existing rows 167-194 down to 168-195.
$ws.Rows("167:167").Insert()

# Populate the newly inserted row 167 with the new record.
$ws.Range("A167").Value = 10
$ws.Range("B167").Value = "Vega Modelo de Temuco"
$ws.Range("C167").Value = "La Araucanía"
$ws.Range("D167").Value = 44995
$ws.Range("E167").Value = 9
$ws.Range("F167").Value = 100112031
$ws.Range("G167").Value = "Poroto verde"
$ws.Range("H167").Value = "Brío"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 200
$ws.Range("K167").Value = 1400
$ws.Range("L167").Value = 1400
$ws.Range("M167").Value = 1400
$ws.Range("N167").Value = "`$/kilo"
$ws.Range("O167").Value = "Región de La Araucanía"
$ws.Range("P167").Value = 1400
$ws.Range("Q167").Value = 1
$ws.Range("R167").Value = "Hortaliza"
